# Apply the Schema.xlsx edits described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Append "/:id" to the three route-path cells in row 6.
$ws.Range("D6").Value = "/api/user/:id"
$ws.Range("H6").Value = "/api/superadmin/changePrivillege/:id"
$ws.Range("L6").Value = "/api/admin/requests/status/:id"

# 2) Remove the "deleteRequest" column (F) from the requests Model/Controller/
#    Routes/Methods block (rows 46-49), shifting the remaining cells left
#    (only within that block, not the whole column).
$ws.Range("F46:G49").Value2 = $ws.Range("G46:H49").Value2
$ws.Range("H46:H49").ClearContents()

# 3) Remove the now-empty trailing formatting-only row 67.
$ws.Rows(67).Delete()

# 4) Widen column H to fit the longer route text.
$ws.Columns("H").ColumnWidth = 34.1640625

# 5) Restore the sheet view / selection state saved with the workbook.
$ws.Application.ActiveWindow.ScrollRow = 38
$ws.Range("F46:G49").Select()
